$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new data point (row) arrived at the top of the table, so the existing
# rows 2..10 (B:G) shift down to rows 3..11 (old row 11's data falls off
# the bottom of the fixed-size table). Walk bottom-up so we never
# overwrite a row before its old value has been copied out.
for ($r = 10; $r -ge 2; $r--) {
    $src = $ws.Range("B$r`:G$r").Value2
    $ws.Range("B$($r+1):G$($r+1)").Value2 = $src
}

# Write the new values for row 2 (the newest data point).
$ws.Range("B2").Value2 = 0.03212158865895828
$ws.Range("C2").Value2 = 0.5134929383467417
$ws.Range("D2").Value2 = 0.5841205420412222
$ws.Range("E2").Value2 = 0.7642777911474481
$ws.Range("F2").Value2 = 0.7845269688049712
$ws.Range("G2").Value2 = 19
